$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking" - Right column corrected from 5 to 4, Wrong column from -1 to -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total" - Right column corrected from 85 to 68, Wrong column from -1 to -2
$ws.Range("B12").Value = 68
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "66 / 112"
